# Shivam Dube (RCB) batting log — add ownTeam/oppTeam columns, refresh rows
# and append two additional match rows (per-match scrape refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("venue", "date", "result", "ownTeam", "oppTeam", "batsman", "totalRuns", "totalBalls", "total4s", "total6s", "sr"),
    @(" Abu Dhabi", " October 28 2020", "Mumbai won by 5 wickets (with 5 balls remaining)", "Royal Challengers Bangalore", "Mumbai Indians", "Shivam Dube ", "2", "6", "0", "0", "33.33"),
    @(" Dubai (DSC)", " September 24 2020", "Kings XI won by 97 runs", "Royal Challengers Bangalore", "Kings XI Punjab", "Shivam Dube ", "12", "12", "0", "1", "100.00"),
    @(" Sharjah", " October 15 2020", "Kings XI won by 8 wickets", "Royal Challengers Bangalore", "Kings XI Punjab", "Shivam Dube ", "23", "19", "0", "2", "121.05"),
    @(" Dubai (DSC)", " September 21 2020", "RCB won by 10 runs", "Royal Challengers Bangalore", "Sunrisers Hyderabad", "Shivam Dube ", "7", "8", "0", "0", "87.50"),
    @(" Abu Dhabi", " November 06 2020", "Sunrisers won by 6 wickets (with 2 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", "Shivam Dube ", "8", "13", "0", "0", "61.53"),
    @(" Dubai (DSC)", " October 05 2020", "Capitals won by 59 runs", "Royal Challengers Bangalore", "Delhi Capitals", "Shivam Dube ", "11", "12", "0", "1", "91.66"),
    @(" Abu Dhabi", " November 02 2020", "Capitals won by 6 wickets (with 6 balls remaining)", "Royal Challengers Bangalore", "Delhi Capitals", "Shivam Dube ", "17", "11", "2", "1", "154.54"),
    @(" Dubai (DSC)", " September 28 2020", "Match tied (RCB won the one-over eliminator)", "Royal Challengers Bangalore", "Mumbai Indians", "Shivam Dube ", "27", "10", "1", "3", "270.00"),
    @(" Dubai (DSC)", " October 10 2020", "RCB won by 37 runs", "Royal Challengers Bangalore", "Chennai Super Kings", "Shivam Dube ", "22", "14", "2", "1", "157.14")
)

# Columns G..K (7..11) hold digit/decimal strings ("8", "61.53", "270.00" …).
# Excel auto-coerces those to numbers on a plain .Value assignment (dropping
# the original text formatting, e.g. "270.00" -> 270), so those columns are
# pre-formatted as Text ("@") before the value is written, matching the
# source sheet where every cell is stored as text.
$numericTextCols = @(7, 8, 9, 10, 11)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $cell = $ws.Cells.Item($r + 1, $c + 1)
        if ($numericTextCols -contains ($c + 1)) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $row[$c]
    }
}
